# EWP happy-path test data for STEPS / PARAMETERS / DATASETS sheets.
$wb = $excel.ActiveWorkbook

$wsTC = $wb.Worksheets.Item("TEST_CASES")
$wsSteps = $wb.Worksheets.Item("STEPS")
$wsParams = $wb.Worksheets.Item("PARAMETERS")
$wsData = $wb.Worksheets.Item("DATASETS")

# A bordered (no font) donor cell from TEST_CASES -- reproduces cellXfs index 2
# (borderId=1) on every cell we touch instead of synthesizing new border defs.
$borderDonor = $wsTC.Range("A3")

# ---- STEPS (sheet2) ----------------------------------------------------
$borderDonor.Copy()
$wsSteps.Range("A1:K9").PasteSpecial(-4122)
$wsSteps.Range("A1:K1").Font.Bold = $true

$wsSteps.Range("A2").Value = "CREATE"
$wsSteps.Range("B2").Value = "owner/path/1"
$wsSteps.Range("E2").Value = 11
$wsSteps.Range("F2").Value = $false
$wsSteps.Range("G2").Value = "action1"
$wsSteps.Range("H2").Value = "result1"

$wsSteps.Range("A3").Value = "C"
$wsSteps.Range("B3").Value = "owner/path/2"
$wsSteps.Range("E3").Value = 12
$wsSteps.Range("F3").Value = 0
$wsSteps.Range("G3").Value = "action2"
$wsSteps.Range("H3").Value = "result2"

$wsSteps.Range("A4").Value = "UPDATE"
$wsSteps.Range("B4").Value = "owner/path/3"
$wsSteps.Range("E4").Value = 13
$wsSteps.Range("F4").Formula = "'0"
$wsSteps.Range("G4").Value = "action3"
$wsSteps.Range("H4").Value = "result3"

$wsSteps.Range("A5").Value = "U"
$wsSteps.Range("B5").Value = "owner/path/4"
$wsSteps.Range("E5").Value = 14
$wsSteps.Range("G5").Value = "action4"
$wsSteps.Range("H5").Value = "result4"

$wsSteps.Range("A6").Value = "DELETE"
$wsSteps.Range("B6").Value = "owner/path/5"
$wsSteps.Range("E6").Value = 15
$wsSteps.Range("G6").Value = "action5"
$wsSteps.Range("H6").Value = "result5"

$wsSteps.Range("A7").Value = "D"
$wsSteps.Range("B7").Value = "owner/path/6"
$wsSteps.Range("E7").Value = 16
$wsSteps.Range("G7").Value = "action6"
$wsSteps.Range("H7").Value = "result6"

$wsSteps.Range("B8").Value = "owner/path/7"
$wsSteps.Range("E8").Value = 17
$wsSteps.Range("G8").Value = "action7"
$wsSteps.Range("H8").Value = "result7"

$wsSteps.Range("A9").Value = "WRONG"
$wsSteps.Range("B9").Value = "owner/path/8"
$wsSteps.Range("E9").Value = 18
$wsSteps.Range("G9").Value = "action8"
$wsSteps.Range("H9").Value = "result8"

$wsSteps.PageSetup.PaperSize = 9
$wsSteps.PageSetup.Orientation = 1

# ---- PARAMETERS (sheet3) ------------------------------------------------
$borderDonor.Copy()
$wsParams.Range("A1:F9").PasteSpecial(-4122)
$wsParams.Range("A1:F1").Font.Bold = $true

$wsParams.Range("A2").Value = "CREATE"
$wsParams.Range("B2").Value = "owner/path/1"
$wsParams.Range("E2").Value = "name1"
$wsParams.Range("F2").Value = "desc1"

$wsParams.Range("A3").Value = "C"
$wsParams.Range("B3").Value = "owner/path/2"
$wsParams.Range("E3").Value = "name2"
$wsParams.Range("F3").Value = "desc2"

$wsParams.Range("A4").Value = "UPDATE"
$wsParams.Range("B4").Value = "owner/path/3"
$wsParams.Range("E4").Value = "name3"
$wsParams.Range("F4").Value = "desc3"

$wsParams.Range("A5").Value = "U"
$wsParams.Range("B5").Value = "owner/path/4"
$wsParams.Range("E5").Value = "name4"
$wsParams.Range("F5").Value = "desc4"

$wsParams.Range("A6").Value = "DELETE"
$wsParams.Range("B6").Value = "owner/path/5"
$wsParams.Range("E6").Value = "name5"
$wsParams.Range("F6").Value = "desc5"

$wsParams.Range("A7").Value = "D"
$wsParams.Range("B7").Value = "owner/path/6"
$wsParams.Range("E7").Value = "name6"
$wsParams.Range("F7").Value = "desc6"

$wsParams.Range("B8").Value = "owner/path/7"
$wsParams.Range("E8").Value = "name7"
$wsParams.Range("F8").Value = "desc7"

$wsParams.Range("A9").Value = "WRONG"
$wsParams.Range("B9").Value = "owner/path/8"
$wsParams.Range("E9").Value = "name8"
$wsParams.Range("F9").Value = "desc8"

$wsParams.PageSetup.PaperSize = 9
$wsParams.PageSetup.Orientation = 1

# ---- DATASETS (sheet4) ---------------------------------------------------
$borderDonor.Copy()
$wsData.Range("A1:I9").PasteSpecial(-4122)
$wsData.Range("A1:I1").Font.Bold = $true

$wsData.Range("A2").Value = "CREATE"
$wsData.Range("B2").Value = "owner/path/1"
$wsData.Range("E2").Value = "name1"
$wsData.Range("G2").Value = "value1"
$wsData.Range("H2").Value = "paramName1"
$wsData.Range("I2").Value = "param/owner/path/1"

$wsData.Range("A3").Value = "C"
$wsData.Range("B3").Value = "owner/path/2"
$wsData.Range("E3").Value = "name2"
$wsData.Range("G3").Value = "value2"
$wsData.Range("H3").Value = "paramName2"
$wsData.Range("I3").Value = "param/owner/path/2"

$wsData.Range("A4").Value = "UPDATE"
$wsData.Range("B4").Value = "owner/path/3"
$wsData.Range("E4").Value = "name3"
$wsData.Range("G4").Value = "value3"
$wsData.Range("H4").Value = "paramName3"
$wsData.Range("I4").Value = "param/owner/path/3"

$wsData.Range("A5").Value = "U"
$wsData.Range("B5").Value = "owner/path/4"
$wsData.Range("E5").Value = "name4"
$wsData.Range("G5").Value = "value4"
$wsData.Range("H5").Value = "paramName4"
$wsData.Range("I5").Value = "param/owner/path/4"

$wsData.Range("A6").Value = "DELETE"
$wsData.Range("B6").Value = "owner/path/5"
$wsData.Range("E6").Value = "name5"
$wsData.Range("G6").Value = "value5"
$wsData.Range("H6").Value = "paramName5"
$wsData.Range("I6").Value = "param/owner/path/5"

$wsData.Range("A7").Value = "D"
$wsData.Range("B7").Value = "owner/path/6"
$wsData.Range("E7").Value = "name6"
$wsData.Range("G7").Value = "value6"
$wsData.Range("H7").Value = "paramName6"
$wsData.Range("I7").Value = "param/owner/path/6"

$wsData.Range("B8").Value = "owner/path/7"
$wsData.Range("E8").Value = "name7"
$wsData.Range("G8").Value = "value7"
$wsData.Range("H8").Value = "paramName7"
$wsData.Range("I8").Value = "param/owner/path/7"

$wsData.Range("A9").Value = "WRONG"
$wsData.Range("B9").Value = "owner/path/8"
$wsData.Range("E9").Value = "name8"
$wsData.Range("G9").Value = "value8"
$wsData.Range("H9").Value = "paramName8"
$wsData.Range("I9").Value = "param/owner/path/8"

$wsData.Columns.Item(9).AutoFit()

$wsData.PageSetup.PaperSize = 9
$wsData.PageSetup.Orientation = 1

# ---- sheet view / selection state ---------------------------------------
$wsParams.Range("E2:E9").Select()
$wsData.Range("C8").Select()
$wsSteps.Range("G17").Select()
